$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 04713e9f-...md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 23:08:50"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
# for 04713e9f-...c710f3c29770016a3250d01c983f2601e4c25fbc.zh-cn.xlf
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-03 23:08:45"
$wsZhCn.Range("K2").Value = "2016-09-03 23:09:05"

# de-de sheet: Correspond Handoff Datetime (H2, shared with Overview!G2) and
# Correspond Handback DateTime (K2) for 04713e9f-...c710f3c29770016a3250d01c983f2601e4c25fbc.de-de.xlf
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 23:08:50"
$wsDeDe.Range("K2").Value = "2016-09-03 23:09:15"
